$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Add a date stamp in C1 (short date format) on the "About" sheet -
# equivalent to the new cell/style introduced by the diff
# (value 44307 = 4/21/2021, numFmtId 14 "short date").
$ws.Range("C1").NumberFormat = "mm-dd-yy"
$ws.Range("C1").Value = Get-Date -Year 2021 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0
